$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the xpath test-data cells so the literal quoting style changes
# from double quotes around the form id to single quotes, matching the
# new test case values used by AreaNameCasesFactory.
$ws.Range("B2").Value = "//form[@id='blog_option_sort_form']//button[span='정확도']"
$ws.Range("B3").Value = "//form[@id='blog_option_sort_form']//button[span='최신순']"
$ws.Range("B4").Value = "//form[@id='blog_option_sort_form']//*[@class='clo_op']"
$ws.Range("B5").Value = "//form[@id='blog_option_sort_form']//*[@class='clo_op']"
$ws.Range("D1").Value = "result"
$ws.Range("D6").Value = "end"

# Move the active selection to B2 (was D6).
$ws.Range("B2").Select()
